$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Foo"
$ws.Range("A2").Value = ""
$ws.Range("A3").Value = "Bar"
$ws.Range("A4").Select()
